$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Showtimes")

# C2 held a date serial (42125, formatted mmm-yy / numFmtId 17). The author
# retyped it as free text ("June 3"), so the cell format flips to Text (@)
# and the value becomes a literal shared string instead of a numeric date.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "June 3"

# Re-measuring the edited rows drops their stale explicit row height so they
# fall back to the sheet default (matches the diff's removal of ht="15.6").
$ws.Rows("1:2").AutoFit()

# Selection moved to E2 in the saved file.
$ws.Range("E2").Select()
